$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.841.72'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').Value = '2.030.78'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.54'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.611'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.44'
$ws.Range('E7').Value = '  +8.42%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0813'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.61'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '2.327.85'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.09'
$ws.Range('E14').Value = '  +4.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.754'
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.24'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').Value = '2.043.05'
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').Value = '37.724.38'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.55'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').Value = '0.0₃0826'
$ws.Range('E21').Value = '  +0.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.22'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').Value = '  -0.71%  '
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.73'
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.14'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.129'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.87'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.28'
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.44'
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.04'
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.51'
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0603'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.28'
$ws.Range('E36').Value = '  +9.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.27'
$ws.Range('E37').Value = '  -1.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.22'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').Value = '1.530.77'
$ws.Range('E40').Value = '  +4.08%  '
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.90'
$ws.Range('E42').Value = '  +2.38%  '
$ws.Range('E43').Value = '  +1.25%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0917'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.79'
$ws.Range('E45').Value = '  -1.55%  '
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.99'
$ws.Range('E47').Value = '  -5.54%  '
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.03'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('D51').Value = '2.218.20'
$ws.Range('E51').Value = '  -0.12%  '
